# Apply scheduled market-price / profit recalculation updates to the Titan_Profits leve sheets.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 50000
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()
$ws.Range("H23").Value = 50000
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()
$ws.Range("H76").Value = 5053577
$ws.Range("I76").Value = 5850890
$ws.Range("K76").Value = 5850890
$ws.Range("M76").Value = -5850575
$ws.Range("H79").Value = 5053577
$ws.Range("I79").Value = 5850890
$ws.Range("K79").Value = 5850890
$ws.Range("M79").Value = -5849798
$ws.Range("H116").Value = 8650515
$ws.Range("I116").Value = 9227016
$ws.Range("J116").Value = 3000
$ws.Range("K116").Value = 9227016
$ws.Range("L116").Value = 3000
$ws.Range("M116").Value = -9223574
$ws.Range("N116").Value = -9884
$ws.Range("H132").Value = 471711.12
$ws.Range("I132").Value = 676288.5600000001
$ws.Range("K132").Value = 2028865.68
$ws.Range("M132").Value = -2026335.68
$ws.Range("H137").Value = 37038612
$ws.Range("I137").Value = 41667764
$ws.Range("J137").Value = 5401
$ws.Range("K137").Value = 125003292
$ws.Range("L137").Value = 16203
$ws.Range("M137").Value = -125000742
$ws.Range("N137").Value = -21303

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14654.448
$ws.Range("I32").Value = 1627.2253
$ws.Range("J32").Value = 146787.72
$ws.Range("K32").Value = 1627.2253
$ws.Range("L32").Value = 146787.72
$ws.Range("M32").Value = -1340.2253
$ws.Range("N32").Value = -147361.72
$ws.Range("H76").Value = 216169
$ws.Range("J76").Value = 216169
$ws.Range("L76").Value = 216169
$ws.Range("N76").Value = -216845
$ws.Range("H79").Value = 216169
$ws.Range("J79").Value = 216169
$ws.Range("L79").Value = 216169
$ws.Range("N79").Value = -218509
$ws.Range("H109").Value = 269877
$ws.Range("J109").Value = 269877
$ws.Range("L109").Value = 269877
$ws.Range("N109").Value = -272651
$ws.Range("H123").Value = 45000
$ws.Range("J123").Value = 45000
$ws.Range("L123").Value = 45000
$ws.Range("N123").Value = -54800

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 950.2222
$ws.Range("I20").Value = 926.43475
$ws.Range("J20").Value = 992.3077
$ws.Range("K20").Value = 926.43475
$ws.Range("L20").Value = 992.3077
$ws.Range("M20").Value = -679.43475
$ws.Range("N20").Value = -1486.3077
$ws.Range("H94").Value = 1043.3529
$ws.Range("I94").Value = 1199.2727
$ws.Range("K94").Value = 1199.2727
$ws.Range("M94").Value = -748.2727
$ws.Range("H105").Value = 3712.9583
$ws.Range("I105").Value = 3473.3333
$ws.Range("K105").Value = 3473.3333
$ws.Range("M105").Value = -1726.3333
$ws.Range("H107").Value = 946.93335
$ws.Range("I107").Value = 698.8
$ws.Range("J107").Value = 1443.2
$ws.Range("K107").Value = 698.8
$ws.Range("L107").Value = 1443.2
$ws.Range("M107").Value = 1221.2
$ws.Range("N107").Value = -5283.2

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2188.9092
$ws.Range("I58").Value = 1629.7059
$ws.Range("K58").Value = 1629.7059
$ws.Range("M58").Value = -1426.7059
$ws.Range("H74").Value = 29125
$ws.Range("J74").Value = 29125
$ws.Range("L74").Value = 29125
$ws.Range("N74").Value = -30873
$ws.Range("H77").Value = 29125
$ws.Range("J77").Value = 29125
$ws.Range("L77").Value = 87375
$ws.Range("N77").Value = -96111
$ws.Range("H88").Value = 30000
$ws.Range("J88").Value = 30000
$ws.Range("L88").Value = 30000
$ws.Range("N88").Value = -30812
$ws.Range("H91").Value = 30000
$ws.Range("J91").Value = 30000
$ws.Range("L91").Value = 30000
$ws.Range("N91").Value = -32808
$ws.Range("H122").Value = 1667.8422
$ws.Range("I122").Value = 944.1539
$ws.Range("J122").Value = 3235.8333
$ws.Range("K122").Value = 2832.4617
$ws.Range("L122").Value = 9707.499899999999
$ws.Range("M122").Value = -382.4616999999998
$ws.Range("N122").Value = -14607.4999
$ws.Range("H136").Value = 2188.9092
$ws.Range("I136").Value = 1629.7059
$ws.Range("K136").Value = 4889.1177
$ws.Range("M136").Value = -2339.1177

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 483.33334
$ws.Range("I92").Value = 325
$ws.Range("K92").Value = 975
$ws.Range("M92").Value = 273
$ws.Range("H132").Value = 1093.1
$ws.Range("I132").Value = 785
$ws.Range("J132").Value = 1328.7059
$ws.Range("K132").Value = 7065
$ws.Range("L132").Value = 11958.3531
$ws.Range("M132").Value = -4535
$ws.Range("N132").Value = -17018.3531

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 6962193
$ws.Range("I11").Value = 9050001
$ws.Range("J11").Value = 2833.3333
$ws.Range("K11").Value = 9050001
$ws.Range("L11").Value = 2833.3333
$ws.Range("M11").Value = -9049862
$ws.Range("N11").Value = -3111.3333
$ws.Range("H70").Value = 7471.875
$ws.Range("I70").Value = 8081.8184
$ws.Range("J70").Value = 6130
$ws.Range("K70").Value = 8081.8184
$ws.Range("L70").Value = 6130
$ws.Range("M70").Value = -7811.8184
$ws.Range("N70").Value = -6670
$ws.Range("H73").Value = 7471.875
$ws.Range("I73").Value = 8081.8184
$ws.Range("K73").Value = 8081.8184
$ws.Range("M73").Value = -7145.8184
$ws.Range("N73").Value = -8002
$ws.Range("H102").Value = 2464.7273
$ws.Range("I102").Value = 2201.5
$ws.Range("K102").Value = 2201.5
$ws.Range("M102").Value = -579.5

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 620.53845
$ws.Range("I93").Value = 696.7
$ws.Range("J93").Value = 366.66666
$ws.Range("K93").Value = 696.7
$ws.Range("L93").Value = 366.66666
$ws.Range("M93").Value = 551.3
$ws.Range("N93").Value = -2862.66666
$ws.Range("H136").Value = 4800.7095
$ws.Range("I136").Value = 2411.625
$ws.Range("J136").Value = 7349.067
$ws.Range("K136").Value = 7234.875
$ws.Range("L136").Value = 22047.201
$ws.Range("M136").Value = -4684.875
$ws.Range("N136").Value = -27147.201

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 12804.667
$ws.Range("J45").Value = 12804.667
$ws.Range("L45").Value = 12804.667
$ws.Range("N45").Value = -13786.667
$ws.Range("H62").Value = 7154740
$ws.Range("I62").Value = 9104305
$ws.Range("J62").Value = 6333.3335
$ws.Range("K62").Value = 9104305
$ws.Range("L62").Value = 6333.3335
$ws.Range("M62").Value = -9103681
$ws.Range("N62").Value = -7581.3335
$ws.Range("H63").Value = 25732.555
$ws.Range("J63").Value = 28324.125
$ws.Range("L63").Value = 28324.125
$ws.Range("N63").Value = -29572.125
$ws.Range("H65").Value = 7154740
$ws.Range("I65").Value = 9104305
$ws.Range("J65").Value = 6333.3335
$ws.Range("K65").Value = 45521525
$ws.Range("L65").Value = 31666.6675
$ws.Range("M65").Value = -45518405
$ws.Range("N65").Value = -37906.6675
$ws.Range("H66").Value = 25732.555
$ws.Range("J66").Value = 28324.125
$ws.Range("L66").Value = 84972.375
$ws.Range("N66").Value = -91212.375
$ws.Range("H82").Value = 35225.75
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 35225.75
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 35225.75
$ws.Range("M82").ClearContents()
$ws.Range("N82").Value = -35991.75
$ws.Range("H85").Value = 35225.75
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 35225.75
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 35225.75
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value = -37877.75
$ws.Range("H132").Value = 13517065
$ws.Range("I132").Value = 29417120
$ws.Range("J132").Value = 2018.5
$ws.Range("K132").Value = 88251360
$ws.Range("L132").Value = 6055.5
$ws.Range("M132").Value = -88248830
$ws.Range("N132").Value = -11115.5
